# [Kadastro App] Yeni kayıt eklendi: 2965
$wb = $excel.ActiveWorkbook

$newRow = @("2965", "2025-09-10", "Erdemli", "1", "ÇAP", "AYHAN KARADAYI (K.Teknisyeni)")
$targetRow = 33

foreach ($sheetName in @("Kayitlar", "Erdemli")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($c = 1; $c -le 6; $c++) {
        $cell = $ws.Cells.Item($targetRow, $c)
        # Force text storage so numeric-looking / date-looking values
        # ("2965", "2025-09-10", "1") are kept as literal text, matching
        # the rest of the column (all cells stored as text).
        $cell.NumberFormat = "@"
        $cell.Value = $newRow[$c - 1]
    }
}
